$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Update existing "想去人数" (want-to-go count) figures
    $ws.Cells.Item(2, 6).Value = 36
    $ws.Cells.Item(3, 6).Value = 80

    # Add a new row (row 5) for the newly announced event, cloning the
    # formatting of the row above (row 4) for column A's numbering cell.
    $ws.Cells.Item(4, 1).Copy($ws.Cells.Item(5, 1))
    $ws.Cells.Item(5, 1).Value = 4

    # Assign the date column as literal text (matching the other rows,
    # which store it as a plain string rather than a date serial), then
    # strip the resulting formatting so the cell ends up with no style
    # index, same as the sibling cells in rows 2-4.
    $ws.Cells.Item(5, 2).NumberFormat = "@"
    $ws.Cells.Item(5, 2).Value = "2024-08-17"
    $ws.Cells.Item(5, 2).ClearFormats()
    $ws.Cells.Item(5, 3).Value = "丽水·AEO纯白礼赞动漫嘉年华"
    $ws.Cells.Item(5, 4).Value = "城北街1001号 爱依·时尚婚宴中心"
    $ws.Cells.Item(5, 5).Value = "2024.08.17 09:00-08.17 16:00"
    $ws.Cells.Item(5, 6).Value = 11
    $ws.Cells.Item(5, 7).Value = 55
    $ws.Cells.Item(5, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86779"
    $ws.Cells.Item(5, 9).Value = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg"
}
